$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.055.04'
$ws.Range("E2").Value = '  -3.85%  '
$ws.Range("D3").Value = '3.142.34'
$ws.Range("E3").Value = '  -3.31%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.47'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.26%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '3.138.51'
$ws.Range("E8").Value = '  -3.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.524'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.75%  '
$ws.Range("E10").Value = '  -6.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.47'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.472'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.36%  '
$ws.Range("E13").Value = '  -8.28%  '
$ws.Range("E14").Value = '  -9.51%  '
$ws.Range("D15").Value = '3.656.75'
$ws.Range("E15").Value = '  -3.39%  '
$ws.Range("D16").Value = '64.074.10'
$ws.Range("E16").Value = '  -3.89%  '
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("D18").Value = '3.138.55'
$ws.Range("E18").Value = '  -3.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.88'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '477.87'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.709'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.54'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.10%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  -5.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.38'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.17'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.73'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.66%  '
$ws.Range("E31").Value = '  -19.22%  '
$ws.Range("E32").Value = '  -6.01%  '
$ws.Range("E33").Value = '  +0.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.08'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.12%  '
$ws.Range("E35").Value = '  -4.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '54.12'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.92'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.19%  '
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = '0.0₃0730'
$ws.Range("E38").Value = '  -8.28%  '
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '460.81'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.94'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -13.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0393'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.78%  '
$ws.Range("E42").Value = '  -8.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.39'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.93%  '
$ws.Range("D44").Value = '2.840.55'
$ws.Range("E44").Value = '  -4.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.264'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -9.86%  '
$ws.Range("E46").Value = '  -10.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.60%  '
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.114'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.89%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.67%  '
